$d = $word.ActiveDocument

$d.Content.Find.Execute("831÷5=166, 1", $true, $false, $false, $false, $false, $true, 1, $false, "159÷8=19, 7", 2) | Out-Null
$d.Content.Find.Execute("555÷4=138, 3", $true, $false, $false, $false, $false, $true, 1, $false, "793÷7=113, 2", 2) | Out-Null
$d.Content.Find.Execute("230÷9=25, 5", $true, $false, $false, $false, $false, $true, 1, $false, "732÷8=91, 4", 2) | Out-Null
$d.Content.Find.Execute("461÷3=153, 2", $true, $false, $false, $false, $false, $true, 1, $false, "857÷2=428, 1", 2) | Out-Null
$d.Content.Find.Execute("551÷4=137, 3", $true, $false, $false, $false, $false, $true, 1, $false, "323÷7=46, 1", 2) | Out-Null
$d.Content.Find.Execute("547÷6=91, 1", $true, $false, $false, $false, $false, $true, 1, $false, "949÷7=135, 4", 2) | Out-Null
$d.Content.Find.Execute("942÷7=134, 4", $true, $false, $false, $false, $false, $true, 1, $false, "324÷6=54, 0", 2) | Out-Null
$d.Content.Find.Execute("816÷6=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "404÷8=50, 4", 2) | Out-Null
$d.Content.Find.Execute("576÷8=72, 0", $true, $false, $false, $false, $false, $true, 1, $false, "240÷9=26, 6", 2) | Out-Null
$d.Content.Find.Execute("181÷3=60, 1", $true, $false, $false, $false, $false, $true, 1, $false, "286÷7=40, 6", 2) | Out-Null
$d.Content.Find.Execute("788÷9=87, 5", $true, $false, $false, $false, $false, $true, 1, $false, "926÷9=102, 8", 2) | Out-Null
$d.Content.Find.Execute("381÷2=190, 1", $true, $false, $false, $false, $false, $true, 1, $false, "345÷2=172, 1", 2) | Out-Null
$d.Content.Find.Execute("719÷2=359, 1", $true, $false, $false, $false, $false, $true, 1, $false, "993÷7=141, 6", 2) | Out-Null
$d.Content.Find.Execute("416÷4=104, 0", $true, $false, $false, $false, $false, $true, 1, $false, "421÷3=140, 1", 2) | Out-Null
$d.Content.Find.Execute("708÷9=78, 6", $true, $false, $false, $false, $false, $true, 1, $false, "566÷4=141, 2", 2) | Out-Null
$d.Content.Find.Execute("161÷7=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "729÷2=364, 1", 2) | Out-Null
$d.Content.Find.Execute("535÷3=178, 1", $true, $false, $false, $false, $false, $true, 1, $false, "187÷8=23, 3", 2) | Out-Null
$d.Content.Find.Execute("943÷9=104, 7", $true, $false, $false, $false, $false, $true, 1, $false, "458÷8=57, 2", 2) | Out-Null
$d.Content.Find.Execute("159÷7=22, 5", $true, $false, $false, $false, $false, $true, 1, $false, "958÷7=136, 6", 2) | Out-Null
$d.Content.Find.Execute("443÷8=55, 3", $true, $false, $false, $false, $false, $true, 1, $false, "536÷8=67, 0", 2) | Out-Null
$d.Content.Find.Execute("838÷8=104, 6", $true, $false, $false, $false, $false, $true, 1, $false, "575÷5=115, 0", 2) | Out-Null
$d.Content.Find.Execute("794÷9=88, 2", $true, $false, $false, $false, $false, $true, 1, $false, "527÷3=175, 2", 2) | Out-Null
$d.Content.Find.Execute("698÷2=349, 0", $true, $false, $false, $false, $false, $true, 1, $false, "470÷3=156, 2", 2) | Out-Null
$d.Content.Find.Execute("725÷3=241, 2", $true, $false, $false, $false, $false, $true, 1, $false, "949÷4=237, 1", 2) | Out-Null
$d.Content.Find.Execute("153÷7=21, 6", $true, $false, $false, $false, $false, $true, 1, $false, "311÷6=51, 5", 2) | Out-Null
